$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Action1")

# B3 previously held "f2", update it to "f3"
$ws.Range("B3").Value = "f3"

# B4 previously held "g4", update it to "g" (same text as A4)
$ws.Range("B4").Value = "g"

# Move the active selection from B4 to B3
$ws.Range("B3").Select()
